# Add a "Save" column (H) to the s_vals sheet, matching the header style
# used by the existing columns (B1:G1) and seeding the data row with 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1: same text style as the other header cells (bold,
# bordered, center/top aligned) -- copy formatting from the adjacent
# header cell G1 so the new column reuses the existing style record
# instead of creating a near-duplicate one.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# New data cell H2: numeric value for the "Save" column.
$ws.Range("H2").Value = 0
